$d = $word.ActiveDocument

# Locate the "Edison Achalma" author paragraph under the title heading
# (the first occurrence in the document, right after "Editar: Editar").
$rng = $d.Content
$found = $rng.Find.Execute("Edison Achalma", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the 'Edison Achalma' author paragraph"
}

# Collapse to the end of the found text, then insert a new paragraph
# (same "Author" style is inherited automatically) containing the
# affiliation line.
$rng.Collapse(0)
$rng.InsertAfter([char]13 + "Escuela Profesional de Economía, Universidad Nacional de San Cristóbal de Huamanga")
